$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert "Interpreter request" row (alphabetically before "Living will", which is row 13) ---
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Value2 = "Interpreter request"

# --- Insert "Motion" row (alphabetically before "Request time off work due to domestic violence",
#     which is at row 15 now that the previous insert shifted it down by one) ---
$ws.Rows.Item(15).Insert()
$ws.Range("A15").Value2 = "Motion"

# Fill in the URL text for both new rows (displayed cell text mirrors the hyperlink URL,
# matching the convention used by every other row in the sheet).
$ws.Range("B13").Value2 = "https://www.illinoislegalaid.org/legal-information/interpreter-request"
$ws.Range("B15").Value2 = "https://www.illinoislegalaid.org/legal-information/motion"

# The engine's row Insert() shifts cell contents/styles but does not shift the worksheet's
# <hyperlinks> ref collection, so rebuild hyperlinks from scratch at their correct
# (post-insert) rows rather than relying on the automatic shift.
$ws.Hyperlinks.Delete()

$hyperlinkTargets = @{
    2  = "https://www.illinoislegalaid.org/legal-information/appearance";
    4  = "https://www.illinoislegalaid.org/legal-information/request-collection-agency-stop-contacting";
    5  = "https://www.illinoislegalaid.org/legal-information/collection-proof-debtor-letter";
    6  = "https://www.illinoislegalaid.org/legal-information/e-filing-exemption-circuit-court";
    7  = "https://www.illinoislegalaid.org/legal-information/e-filing-exemption-appellate-court";
    8  = "https://www.illinoislegalaid.org/legal-information/e-filing-exemption-supreme-court";
    9  = "https://www.illinoislegalaid.org/legal-information/end-illegal-lockout-demand";
    11 = "https://www.illinoislegalaid.org/legal-information/fee-waiver";
    12 = "https://www.illinoislegalaid.org/legal-information/housing-discrimination-complaint-idhr";
    13 = "https://www.illinoislegalaid.org/legal-information/interpreter-request";
    15 = "https://www.illinoislegalaid.org/legal-information/motion";
    16 = "https://www.illinoislegalaid.org/legal-information/request-time-work-due-domestic-abuse-letter";
    17 = "https://www.illinoislegalaid.org/legal-information/respond-lawsuit";
    18 = "https://www.illinoislegalaid.org/legal-information/security-deposit-demand-letter";
    19 = "https://www.illinoislegalaid.org/legal-information/stop-wage-assignment-letter";
    20 = "https://www.illinoislegalaid.org/legal-information/voluntary-acknowledgment-parentage-vap";
}

# Recreate hyperlinks in the same relative order as the original workbook (by row number)
# so the relationship ids come out in the same sequence, then restore the "Hyperlink" cell
# style that Hyperlinks.Add() mutates away from the shared style used by the rest of column B.
$orderedRows = @(2, 11, 5, 4, 9, 18, 12, 19, 16, 6, 7, 8, 17, 20, 13, 15)
foreach ($r in $orderedRows) {
    $target = $hyperlinkTargets[$r]
    $cell = $ws.Range("B$r")
    $ws.Hyperlinks.Add($cell, $target)
    $cell.Style = $ws.Range("B2").Style
}
